$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to remain text (matches original inlineStr typing)
# so that numeric-looking strings (e.g. "33.30", "8.524") are not
# auto-converted into numbers by Excel when assigned below.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.610.95'
$ws.Range('E2').Value = '  +1.91%  '
$ws.Range('D3').Value = '1.890.51'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').Value = '244.51'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').Value = '0.4952'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '0.2948'
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('D9').Value = '0.06793'
$ws.Range('E9').Value = '  +2.78%  '
$ws.Range('D10').Value = '1.892.32'
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('D11').Value = '17.01'
$ws.Range('E11').Value = '  +1.59%  '
$ws.Range('D12').Value = '0.07298'
$ws.Range('E12').Value = '  +1.79%  '
$ws.Range('D13').Value = '90.55'
$ws.Range('E13').Value = '  +5.25%  '
$ws.Range('D14').Value = '5.044'
$ws.Range('E14').Value = '  +4.18%  '
$ws.Range('D15').Value = '0.6701'
$ws.Range('E15').Value = '  +1.17%  '
$ws.Range('D16').Value = '30.614.31'
$ws.Range('E16').Value = '  +2.04%  '
$ws.Range('D17').Value = '0.000007914'
$ws.Range('E17').Value = '  +0.79%  '
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').Value = '13.17'
$ws.Range('E19').Value = '  +3.63%  '
$ws.Range('D20').Value = '2.141.12'
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('E21').Value = '  +0.56%  '
$ws.Range('D22').Value = '4.845'
$ws.Range('E22').Value = '  +1.89%  '
$ws.Range('D23').Value = '174.72'
$ws.Range('E23').Value = '  +29.97%  '
$ws.Range('D24').Value = '6.043'
$ws.Range('E24').Value = '  +7.89%  '
$ws.Range('D25').Value = '9.243'
$ws.Range('E25').Value = '  +1.56%  '
$ws.Range('D26').Value = '155.26'
$ws.Range('E26').Value = '  +3.14%  '
$ws.Range('D27').Value = '18.53'
$ws.Range('E27').Value = '  +10.82%  '
$ws.Range('D28').Value = '1.920'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').Value = '1.389'
$ws.Range('E29').Value = '  +1.37%  '
$ws.Range('D30').Value = '4.321'
$ws.Range('E30').Value = '  +3.81%  '
$ws.Range('D31').Value = '0.08901'
$ws.Range('E31').Value = '  +2.43%  '
$ws.Range('D32').Value = '4.012'
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('D33').Value = '0.05224'
$ws.Range('E33').Value = '  +4.45%  '
$ws.Range('D34').Value = '0.7376'
$ws.Range('E34').Value = '  +4.48%  '
$ws.Range('D35').Value = '1.130'
$ws.Range('E35').Value = '  +2.77%  '
$ws.Range('D36').Value = '2.674'
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('D37').Value = '0.01865'
$ws.Range('E37').Value = '  +10.67%  '
$ws.Range('D38').Value = '2.701'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').Value = '2.159'
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('D40').Value = '0.9352'
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('D41').Value = '0.4344'
$ws.Range('E41').Value = '  +3.78%  '
$ws.Range('D42').Value = '105.77'
$ws.Range('E42').Value = '  +3.99%  '
$ws.Range('D43').Value = '5.801'
$ws.Range('E43').Value = '  -2.65%  '
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('D45').Value = '7.638'
$ws.Range('E45').Value = '  +3.07%  '
$ws.Range('D46').Value = '0.1356'
$ws.Range('E46').Value = '  +7.89%  '
$ws.Range('D47').Value = '0.05828'
$ws.Range('E47').Value = '  +2.82%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').Value = '33.30'
$ws.Range('E48').Value = '  +2.63%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '8.524'
$ws.Range('E49').Value = '  +5.80%  '
$ws.Range('D50').Value = '0.3865'
$ws.Range('E50').Value = '  +4.50%  '
$ws.Range('D51').Value = '1.377'
$ws.Range('E51').Value = '  +3.10%  '

# Restore default (unstyled) appearance for the data range so the
# cells keep looking identical to how they did before (no borders,
# bold, etc. were ever applied to this range).
$ws.Range("D2:E51").Style = "Normal"
